# Auto-generated Excel COM-interop script.
# Applies the numeric cell updates described by the target diff: value replacements,
# a couple of new-cell insertions, and two cell deletions (WVR N64/N67).
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2870.7646
$ws.Range("H33").Value = 272.73334
$ws.Range("I33").Value = 331.2857
$ws.Range("J33").Value = 221.5
$ws.Range("K33").Value = 331.2857
$ws.Range("L33").Value = 221.5
$ws.Range("M33").Value = -102.2857
$ws.Range("N33").Value = -679.5
$ws.Range("H38").Value = 514.5
$ws.Range("I38").Value = 174.25
$ws.Range("J38").Value = 968.1667
$ws.Range("K38").Value = 522.75
$ws.Range("L38").Value = 2904.5001
$ws.Range("M38").Value = -150.75
$ws.Range("N38").Value = -3648.5001
$ws.Range("H87").Value = 29288.5
$ws.Range("J87").Value = 29288.5
$ws.Range("L87").Value = 29288.5
$ws.Range("N87").Value = -31784.5
$ws.Range("H90").Value = 29288.5
$ws.Range("J90").Value = 29288.5
$ws.Range("L90").Value = 87865.5
$ws.Range("N90").Value = -100345.5
$ws.Range("H129").Value = 914.9530999999999
$ws.Range("I129").Value = 295.14285
$ws.Range("J129").Value = 991.0702
$ws.Range("K129").Value = 885.4285500000001
$ws.Range("L129").Value = 2973.2106
$ws.Range("M129").Value = 4114.571449999999
$ws.Range("N129").Value = -12973.2106
$ws.Range("H137").Value = 567604.9
$ws.Range("I137").Value = 2835.4443
$ws.Range("J137").Value = 1044129.1
$ws.Range("K137").Value = 8506.332900000001
$ws.Range("L137").Value = 3132387.3
$ws.Range("M137").Value = -5956.332900000001
$ws.Range("N137").Value = -3137487.3
$ws.Range("H138").Value = 4946.2705
$ws.Range("I138").Value = 1015.7857
$ws.Range("J138").Value = 7338.7393
$ws.Range("K138").Value = 3047.3571
$ws.Range("L138").Value = 22016.2179
$ws.Range("M138").Value = 2092.6429
$ws.Range("N138").Value = -32296.2179

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15933.162
$ws.Range("I32").Value = 16500.572
$ws.Range("K32").Value = 16500.572
$ws.Range("M32").Value = -16213.572
$ws.Range("H61").Value = 6113.479
$ws.Range("I61").Value = 2336.2778
$ws.Range("J61").Value = 17445.084
$ws.Range("K61").Value = 2336.2778
$ws.Range("L61").Value = 17445.084
$ws.Range("M61").Value = -2124.2778
$ws.Range("N61").Value = -17869.084
$ws.Range("H136").Value = 6113.479
$ws.Range("I136").Value = 2336.2778
$ws.Range("J136").Value = 17445.084
$ws.Range("K136").Value = 7008.8334
$ws.Range("L136").Value = 52335.25199999999
$ws.Range("M136").Value = -4458.8334
$ws.Range("N136").Value = -57435.25199999999

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 48291.47
$ws.Range("J132").Value = 48291.47
$ws.Range("L132").Value = 48291.47
$ws.Range("N132").Value = -58411.47
$ws.Range("H134").Value = 47681.047
$ws.Range("I134").Value = 2423.7896
$ws.Range("J134").Value = 334310.34
$ws.Range("K134").Value = 7271.3688
$ws.Range("L134").Value = 1002931.02
$ws.Range("M134").Value = -4736.3688
$ws.Range("N134").Value = -1008001.02
$ws.Range("H139").Value = 40000
$ws.Range("J139").Value = 40000
$ws.Range("L139").Value = 40000
$ws.Range("N139").Value = -50280
$ws.Range("H141").Value = 39865
$ws.Range("J141").Value = 39865
$ws.Range("L141").Value = 39865
$ws.Range("N141").Value = -50225

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 19000
$ws.Range("J36").Value = 19000
$ws.Range("L36").Value = 19000
$ws.Range("N36").Value = -19776
$ws.Range("H40").Value = 19000
$ws.Range("J40").Value = 19000
$ws.Range("L40").Value = 19000
$ws.Range("N40").Value = -19320
$ws.Range("H58").Value = 3959170.2
$ws.Range("I58").Value = 5349818
$ws.Range("J58").Value = 19002.334
$ws.Range("K58").Value = 5349818
$ws.Range("L58").Value = 19002.334
$ws.Range("M58").Value = -5349615
$ws.Range("N58").Value = -19408.334
$ws.Range("H107").Value = 879.56525
$ws.Range("J107").Value = 1002
$ws.Range("L107").Value = 1002
$ws.Range("N107").Value = -4842
$ws.Range("H136").Value = 3959170.2
$ws.Range("I136").Value = 5349818
$ws.Range("J136").Value = 19002.334
$ws.Range("K136").Value = 16049454
$ws.Range("L136").Value = 57007.00199999999
$ws.Range("M136").Value = -16046904
$ws.Range("N136").Value = -62107.00199999999

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 723
$ws.Range("I52").Value = 372
$ws.Range("J52").Value = 810.75
$ws.Range("K52").Value = 1116
$ws.Range("L52").Value = 2432.25
$ws.Range("M52").Value = -850
$ws.Range("N52").Value = -2964.25
$ws.Range("H68").Value = 1210.8673
$ws.Range("I68").Value = 1229.1224
$ws.Range("J68").Value = 1192.6123
$ws.Range("K68").Value = 3687.3672
$ws.Range("L68").Value = 3577.8369
$ws.Range("M68").Value = -2876.3672
$ws.Range("N68").Value = -5199.8369
$ws.Range("H69").Value = 100002620
$ws.Range("I69").Value = 1100
$ws.Range("J69").Value = 125003000
$ws.Range("K69").Value = 3300
$ws.Range("L69").Value = 375009000
$ws.Range("M69").Value = -2489
$ws.Range("N69").Value = -375010622
$ws.Range("H71").Value = 1210.8673
$ws.Range("I71").Value = 1229.1224
$ws.Range("J71").Value = 1192.6123
$ws.Range("K71").Value = 11062.1016
$ws.Range("L71").Value = 10733.5107
$ws.Range("M71").Value = -7006.1016
$ws.Range("N71").Value = -18845.5107
$ws.Range("H72").Value = 100002620
$ws.Range("I72").Value = 1100
$ws.Range("J72").Value = 125003000
$ws.Range("K72").Value = 9900
$ws.Range("L72").Value = 1125027000
$ws.Range("M72").Value = -5844
$ws.Range("N72").Value = -1125035112
$ws.Range("H134").Value = 4132.5654

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1250
$ws.Range("I22").Value = 1250
$ws.Range("K22").Value = 1250
$ws.Range("M22").Value = -955
$ws.Range("H27").Value = 1250
$ws.Range("I27").Value = 1250
$ws.Range("K27").Value = 1250
$ws.Range("M27").Value = -1143

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H101").Value = 26097
$ws.Range("J101").Value = 26097
$ws.Range("L101").Value = 26097
$ws.Range("N101").Value = -32587.0
